$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage so numeric-looking strings (e.g. "1.001") are not
# reinterpreted as numbers by Excel's input parser.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.584.54"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "1.883.04"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "246.17"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "0.4740"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "0.2893"
$ws.Range("E8").Value = "  -1.44%  "
$ws.Range("D9").Value = "0.06531"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "22.30"
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("D11").Value = "0.7612"
$ws.Range("E11").Value = "  +2.51%  "
$ws.Range("D12").Value = "99.51"
$ws.Range("E12").Value = "  +2.45%  "
$ws.Range("D13").Value = "0.07823"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").Value = "1.882.35"
$ws.Range("E14").Value = "  -0.76%  "
$ws.Range("D15").Value = "5.229"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "283.51"
$ws.Range("E16").Value = "  -1.08%  "
$ws.Range("D17").Value = "30.567.69"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "13.16"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").Value = "0.000007515"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "2.127.51"
$ws.Range("D22").Value = "5.348"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "1.001"
$ws.Range("D24").Value = "6.428"
$ws.Range("E24").Value = "  +2.15%  "
$ws.Range("D25").Value = "9.164"
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("D26").Value = "164.02"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").Value = "19.00"
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").Value = "1.904"
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("D29").Value = "0.09748"
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("D30").Value = "1.327"
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("D31").Value = "1.502"
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("D32").Value = "4.247"
$ws.Range("E32").Value = "  -1.60%  "
$ws.Range("D33").Value = "4.180"
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("D34").Value = "0.04840"
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("D35").Value = "1.131"
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("D36").Value = "0.6976"
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("D37").Value = "2.772"
$ws.Range("E37").Value = "  +2.36%  "
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").Value = "2.874"
$ws.Range("E39").Value = "  +1.19%  "
$ws.Range("D40").Value = "6.302"
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("D41").Value = "75.35"
$ws.Range("E41").Value = "  -0.98%  "
$ws.Range("D42").Value = "1.972"
$ws.Range("E42").Value = "  -2.28%  "
$ws.Range("D43").Value = "0.4245"
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "0.8388"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").Value = "9.987"
$ws.Range("E46").Value = "  +3.89%  "
$ws.Range("D47").Value = "101.38"
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("D48").Value = "7.009"
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("D49").Value = "35.24"
$ws.Range("E49").Value = "  -0.67%  "
$ws.Range("D50").Value = "0.05771"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("D51").Value = "0.3953"
$ws.Range("E51").Value = "  -0.66%  "

# Restore the default (unstyled) cell style so no stray style index is
# left attached to these cells, matching the original formatting.
$ws.Range("D2:E51").Style = "Normal"
